$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Insert the new "Users" sheet as the very first sheet in the workbook
# (sheet name "Users", placed before the current first sheet "Category").
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "Users"

# ---------------------------------------------------------------------------
# Fill in the login credentials table.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "user1"
$ws.Range("B2").Value = "pass1"
$ws.Range("A3").Value = "user2"
$ws.Range("B3").Value = "pass2"
$ws.Range("A4").Value = "user3"
$ws.Range("B4").Value = "pass3"

# ---------------------------------------------------------------------------
# Style the data rows (A2:B4) with vertical-centered, wrapped text.
# Build the combined alignment format once on a scratch cell, then copy the
# format across with PasteSpecial so only a single new style entry is
# produced (matches the authored workbook's cellXfs count="2").
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4108

$dataRange = $ws.Range("A2:B4")
$scratch.Copy()
$dataRange.PasteSpecial(-4122)
$scratch.EntireRow.Delete()

# ---------------------------------------------------------------------------
# Match the saved selection/active-cell state of the authored sheet.
# ---------------------------------------------------------------------------
[void]$ws.Range("B22").Select()
